$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 8; existing rows 8-15 shift down to 9-16,
# carrying their formatting (incl. the date style on column D) with them.
$ws.Rows("8:8").Insert()

# Populate the newly inserted row 8 with the new weekly record.
$ws.Range("A8").Value = 9
$ws.Range("B8").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C8").Value = "Metropolitana"
$ws.Range("D8").Value = 44978
$ws.Range("E8").Value = 13
$ws.Range("F8").Value = "Fruta"
$ws.Range("G8").Value = 100101
$ws.Range("H8").Value = "Berries"
$ws.Range("I8").Value = 100101008
$ws.Range("J8").Value = "Mora"
$ws.Range("K8").Value = "Sin especificar"
$ws.Range("L8").Value = "Primera"
$ws.Range("M8").Value = 500
$ws.Range("N8").Value = 3000
$ws.Range("O8").Value = 3000
$ws.Range("P8").Value = 3000
$ws.Range("Q8").Value = "$/bandeja 2 kilos"
$ws.Range("R8").Value = "Provincia de Curicó"
$ws.Range("S8").Value = 1500
$ws.Range("T8").Value = 2
